$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "71×24=1704" "20×18=360"
Replace-Text "60×63=3780" "58×46=2668"
Replace-Text "82×54=4428" "17×13=221"
Replace-Text "19×22=418" "26×17=442"
Replace-Text "71×62=4402" "64×46=2944"
Replace-Text "26×65=1690" "65×86=5590"
Replace-Text "86×80=6880" "31×90=2790"
Replace-Text "87×88=7656" "65×84=5460"
Replace-Text "44×13=572" "57×60=3420"
Replace-Text "77×60=4620" "34×66=2244"
Replace-Text "63×15=945" "81×48=3888"
Replace-Text "64×57=3648" "50×90=4500"
Replace-Text "68×11=748" "52×11=572"
Replace-Text "24×91=2184" "48×45=2160"
Replace-Text "22×81=1782" "29×89=2581"
Replace-Text "42×27=1134" "98×83=8134"
Replace-Text "36×86=3096" "86×81=6966"
Replace-Text "45×26=1170" "11×83=913"
Replace-Text "64×33=2112" "21×73=1533"
Replace-Text "21×54=1134" "82×35=2870"
Replace-Text "17×42=714" "55×30=1650"
Replace-Text "62×68=4216" "32×60=1920"
Replace-Text "24×45=1080" "75×33=2475"
Replace-Text "86×45=3870" "88×55=4840"
Replace-Text "92×13=1196" "64×87=5568"

Write-Host "Done replacing"
